$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 626.75
$ws.Range("I9").Value = 692.8
$ws.Range("J9").Value = 516.6667
$ws.Range("K9").Value = 692.8
$ws.Range("L9").Value = 516.6667
$ws.Range("M9").Value = -523.8
$ws.Range("N9").Value = -854.6667
# Row 15
$ws.Range("H15").Value = 1236.5571
$ws.Range("I15").Value = 1236.5571
$ws.Range("K15").Value = 3709.6713
$ws.Range("M15").Value = -3540.6713
# Row 28
$ws.Range("H28").Value = 809.6667
$ws.Range("I28").Value = 443
$ws.Range("J28").Value = 2276.3333
$ws.Range("K28").Value = 443
$ws.Range("L28").Value = 2276.3333
$ws.Range("M28").Value = 42
$ws.Range("N28").Value = -3246.3333
# Row 62
$ws.Range("H62").Value = 3878.5557
$ws.Range("I62").Value = 3311.3333
$ws.Range("K62").Value = 3311.3333
$ws.Range("M62").Value = -2687.3333
# Row 65
$ws.Range("H65").Value = 3878.5557
$ws.Range("I65").Value = 3311.3333
$ws.Range("K65").Value = 16556.6665
$ws.Range("M65").Value = -13436.6665
# Row 80
$ws.Range("H80").Value = 8998.154
$ws.Range("I80").Value = 500
$ws.Range("J80").Value = 12775.111
$ws.Range("K80").Value = 1500
$ws.Range("L80").Value = 38325.333
$ws.Range("M80").Value = -502
$ws.Range("N80").Value = -40321.333
# Row 83
$ws.Range("H83").Value = 8998.154
$ws.Range("I83").Value = 500
$ws.Range("J83").Value = 12775.111
$ws.Range("K83").Value = 4500
$ws.Range("L83").Value = 114975.999
$ws.Range("M83").Value = 492
$ws.Range("N83").Value = -124959.999
# Row 107
$ws.Range("H107").Value = 1365
$ws.Range("I107").Value = 961.1667
$ws.Range("K107").Value = 961.1667
$ws.Range("M107").Value = 958.8333
# Row 112
$ws.Range("H112").Value = 4953.1875
$ws.Range("I112").Value = 1374.5
$ws.Range("J112").Value = 5464.4287
$ws.Range("K112").Value = 4123.5
$ws.Range("L112").Value = 16393.2861
$ws.Range("M112").Value = -3015.5
$ws.Range("N112").Value = -18609.2861
# Row 132
$ws.Range("H132").Value = 2930742.5
$ws.Range("I132").Value = 2930742.5
$ws.Range("K132").Value = 8792227.5
$ws.Range("M132").Value = -8789697.5

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 211167.94
$ws.Range("I74").Value = 231460.42
$ws.Range("K74").Value = 231460.42
$ws.Range("M74").Value = -230586.42
# Row 77
$ws.Range("H77").Value = 211167.94
$ws.Range("I77").Value = 231460.42
$ws.Range("K77").Value = 1157302.1
$ws.Range("M77").Value = -1152934.1
# Row 102
$ws.Range("H102").Value = 3624.4443
$ws.Range("I102").Value = 3813.625
$ws.Range("K102").Value = 3813.625
$ws.Range("M102").Value = -2191.625
# Row 122
$ws.Range("H122").Value = 2005.3846
$ws.Range("J122").Value = 1500
$ws.Range("L122").Value = 4500
$ws.Range("N122").Value = -9400

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2750.35
$ws.Range("I105").Value = 2505.3333
$ws.Range("K105").Value = 2505.3333
$ws.Range("M105").Value = -758.3332999999998

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 166.66667
$ws.Range("J4").Value = 500
$ws.Range("L4").Value = 500
$ws.Range("N4").Value = -724
# Row 6
$ws.Range("H6").Value = 170934.61
$ws.Range("I6").Value = 2000
$ws.Range("J6").Value = 211478.92
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 211478.92
$ws.Range("M6").Value = -1887
$ws.Range("N6").Value = -211704.92
# Row 14
$ws.Range("H14").Value = 4665
$ws.Range("I14").Value = 3995
$ws.Range("K14").Value = 3995
$ws.Range("M14").Value = -3825
# Row 16
$ws.Range("H16").Value = 2427.125
$ws.Range("I16").Value = 1669.25
$ws.Range("K16").Value = 1669.25
$ws.Range("M16").Value = -1382.25
# Row 94
$ws.Range("H94").Value = 1150.7916
$ws.Range("I94").Value = 867.9
$ws.Range("J94").Value = 1352.8572
$ws.Range("K94").Value = 867.9
$ws.Range("L94").Value = 1352.8572
$ws.Range("M94").Value = -416.9
$ws.Range("N94").Value = -2254.8572
# Row 113
$ws.Range("H113").Value = 2427.125
$ws.Range("I113").Value = 1669.25
$ws.Range("K113").Value = 1669.25
$ws.Range("M113").Value = 500.75

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 50
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
# Row 53
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
# Row 56
$ws.Range("H56").Value = 6284.857
$ws.Range("I56").Value = 6284.857
$ws.Range("K56").Value = 6284.857
$ws.Range("M56").Value = -5754.857
# Row 107
$ws.Range("H107").Value = 2442
$ws.Range("J107").Value = 872.75
$ws.Range("L107").Value = 2618.25
$ws.Range("N107").Value = -6458.25

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 182.95238
$ws.Range("I2").Value = 151.5
$ws.Range("J2").Value = 211.54546
$ws.Range("K2").Value = 151.5
$ws.Range("L2").Value = 211.54546
$ws.Range("M2").Value = -38.5
$ws.Range("N2").Value = -437.54546
# Row 46
$ws.Range("H46").Value = 16332.5
# Row 48
$ws.Range("H48").Value = 38705.727
$ws.Range("J48").Value = 39977.3
$ws.Range("L48").Value = 39977.3
$ws.Range("N48").Value = -40947.3
# Row 126
$ws.Range("H126").Value = 2835.5454
$ws.Range("I126").Value = 2023.875
$ws.Range("K126").Value = 6071.625
$ws.Range("M126").Value = -3601.625

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 12000
$ws.Range("J3").Value = 12000
$ws.Range("L3").Value = 12000
$ws.Range("N3").Value = -12224
# Row 7
$ws.Range("H7").Value = 2900.4
$ws.Range("I7").Value = 2853.2632
$ws.Range("K7").Value = 2853.2632
$ws.Range("M7").Value = -2741.2632
# Row 9
$ws.Range("H9").Value = 950
$ws.Range("I9").Value = 900
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 900
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = -676
$ws.Range("N9").Value = -1448
# Row 11
$ws.Range("H11").Value = 39500
$ws.Range("J11").Value = 39500
$ws.Range("L11").Value = 39500
$ws.Range("N11").Value = -39780
# Row 13
$ws.Range("H13").Value = 12000
$ws.Range("J13").Value = 12000
$ws.Range("L13").Value = 12000
$ws.Range("N13").Value = -12280
# Row 14
$ws.Range("H14").Value = 14004.5
$ws.Range("J14").Value = 14004.5
$ws.Range("L14").Value = 14004.5
$ws.Range("N14").Value = -14348.5
# Row 15
$ws.Range("H15").Value = 12000
$ws.Range("J15").Value = 12000
$ws.Range("L15").Value = 12000
$ws.Range("N15").Value = -12340
# Row 19
$ws.Range("H19").Value = 163.33333
$ws.Range("I19").Value = 163.33333
$ws.Range("K19").Value = 163.33333
$ws.Range("M19").Value = 6.666670000000011
# Row 40
$ws.Range("H40").Value = 1884.1111
$ws.Range("I40").Value = 1764.2693
$ws.Range("K40").Value = 1764.2693
$ws.Range("M40").Value = -1628.2693
# Row 46
$ws.Range("H46").Value = 4251.1333
$ws.Range("I46").Value = 900.3333
$ws.Range("J46").Value = 5088.8335
$ws.Range("K46").Value = 900.3333
$ws.Range("L46").Value = 5088.8335
$ws.Range("M46").Value = -712.3333
$ws.Range("N46").Value = -5464.8335
# Row 55
$ws.Range("H55").Value = 670.65216
$ws.Range("J55").Value = 1626.125
$ws.Range("L55").Value = 1626.125
$ws.Range("N55").Value = -1972.125
# Row 61
$ws.Range("H61").Value = 1196.75
$ws.Range("I61").Value = 1196.75
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1196.75
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -994.75
$ws.Range("N61").ClearContents()
# Row 82
$ws.Range("H82").Value = 1603.6923
$ws.Range("I82").Value = 1478.2858
$ws.Range("J82").Value = 1750
$ws.Range("K82").Value = 1478.2858
$ws.Range("L82").Value = 1750
$ws.Range("M82").Value = -1117.2858
$ws.Range("N82").Value = -2472
# Row 85
$ws.Range("H85").Value = 1603.6923
$ws.Range("I85").Value = 1478.2858
$ws.Range("J85").Value = 1750
$ws.Range("K85").Value = 1478.2858
$ws.Range("L85").Value = 1750
$ws.Range("M85").Value = -230.2858000000001
$ws.Range("N85").Value = -4246
# Row 93
$ws.Range("H93").Value = 2154.889
$ws.Range("I93").Value = 1899.2
$ws.Range("J93").Value = 2474.5
$ws.Range("K93").Value = 1899.2
$ws.Range("L93").Value = 2474.5
$ws.Range("M93").Value = -651.2
$ws.Range("N93").Value = -4970.5
# Row 113
$ws.Range("H113").Value = 1196.75
$ws.Range("I113").Value = 1196.75
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1196.75
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 973.25
$ws.Range("N113").ClearContents()
# Row 126
$ws.Range("H126").Value = 2900.4
$ws.Range("I126").Value = 2853.2632
$ws.Range("K126").Value = 8559.7896
$ws.Range("M126").Value = -6089.7896
# Row 132
$ws.Range("H132").Value = 2102.4666
$ws.Range("J132").Value = 3429
$ws.Range("L132").Value = 10287
$ws.Range("N132").Value = -15347
# Row 136
$ws.Range("H136").Value = 3081.889
$ws.Range("I136").Value = 2983.4443
$ws.Range("J136").Value = 3377.2222
$ws.Range("K136").Value = 8950.332900000001
$ws.Range("L136").Value = 10131.6666
$ws.Range("M136").Value = -6400.332900000001
$ws.Range("N136").Value = -15231.6666

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 1602.6
$ws.Range("I100").Value = 1374.6666
$ws.Range("J100").Value = 1944.5
$ws.Range("K100").Value = 2749.3332
$ws.Range("L100").Value = 3889
$ws.Range("M100").Value = -2208.3332
$ws.Range("N100").Value = -4971
# Row 122
$ws.Range("H122").Value = 155806.06
$ws.Range("I122").Value = 231509.2
$ws.Range("K122").Value = 694527.6000000001
$ws.Range("M122").Value = -692077.6000000001
